# Update gh-pages to output generated at 456a3b4
# Refresh the scraped "想去人数" (interest-count) figures in column F on
# the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 13318
$ws1.Range("F6").Value  = 1004
$ws1.Range("F8").Value  = 1718
$ws1.Range("F12").Value = 39
$ws1.Range("F14").Value = 13306
$ws1.Range("F16").Value = 581
$ws1.Range("F17").Value = 8880
$ws1.Range("F18").Value = 1
$ws1.Range("F19").Value = 7958

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 13318
$ws4.Range("F7").Value  = 1004
$ws4.Range("F9").Value  = 1718
$ws4.Range("F13").Value = 39
$ws4.Range("F15").Value = 13306
$ws4.Range("F17").Value = 581
$ws4.Range("F18").Value = 8880
$ws4.Range("F19").Value = 1
$ws4.Range("F20").Value = 7958
